$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'243.38"
$ws.Range('D3').Value = "'23.04"
$ws.Range('D4').Value = "'5.407"
$ws.Range('D6').Value = "'3.422"
$ws.Range('D7').Value = "'6.498"
$ws.Range('D8').Value = "'0.8114"
$ws.Range('D9').Value = "'0.9255"
$ws.Range('B10').Value = "'One"
$ws.Range('C10').Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range('D10').Value = "'0.01115"
$ws.Range('E10').Value = "'9OneONEBestin24h"
$ws.Range('B11').Value = "'WazirX"
$ws.Range('C11').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D11').Value = "'0.1439"
$ws.Range('E11').Value = "'10WazirXWRX"
$ws.Range('B12').Value = "'MandalaExchangeToken"
$ws.Range('C12').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D12').Value = "'0.07428"
$ws.Range('E12').Value = "'11MandalaExchangeTokenMDX"
$ws.Range('B13').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('C13').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D13').Value = "'0.03318"
$ws.Range('E13').Value = "'12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range('B14').Value = "'BitrueCoin"
$ws.Range('C14').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D14').Value = "'0.03097"
$ws.Range('E14').Value = "'13BitrueCoinBTR"
$ws.Range('B15').Value = "'BitMartToken"
$ws.Range('C15').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D15').Value = "'0.09355"
$ws.Range('E15').Value = "'14BitMartTokenBMX"
$ws.Range('B16').Value = "'MCDex"
$ws.Range('C16').Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range('D16').Value = "'3.850"
$ws.Range('E16').Value = "'15MCDexMCB"
$ws.Range('B17').Value = "'BitForexToken"
$ws.Range('C17').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D17').Value = "'0.001571"
$ws.Range('E17').Value = "'16BitForexTokenBF"
$ws.Range('B18').Value = "'CoinExToken"
$ws.Range('C18').Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range('D18').Value = "'0.04726"
$ws.Range('E18').Value = "'17CoinExTokenCET"
$ws.Range('D19').Value = "'0.005852"
$ws.Range('D20').Value = "'0.001270"
$ws.Range('E20').Value = "'19BitKanKAN"
$ws.Range('D21').Value = "'0.004879"
$ws.Range('D23').Value = "'3.572"
$ws.Range('D25').Value = "'0.3237"
$ws.Range('D26').Value = "'0.1331"
$ws.Range('D27').Value = "'0.0002340"
$ws.Range('D40').Value = "'0.03954"
$ws.Range('D41').Value = "'0.006472"
$ws.Range('B42').Value = "'CEJI"
$ws.Range('C42').Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range('D42').Value = "'0.004001"
$ws.Range('E42').Value = "'41CEJICEJI"
$ws.Range('B43').Value = "'BKEXToken"
$ws.Range('C43').Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range('D43').Value = "'0.1076"
$ws.Range('E43').Value = "'42BKEXTokenBKK"
$ws.Range('D44').Value = "'0.008922"
$ws.Range('D45').Value = "'0.00005192"
$ws.Range('D47').Value = "'0.7112"
$ws.Range('D48').Value = "'0.002274"
$ws.Range('D50').Value = "'0.0002000"

Write-Output "Applied 66 cell updates"
